$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new note about generating the three stacks to E13, preserving the
# existing wrap-text style already applied to that cell.
$note = "To generate the three stacks needed for the three pegs, use this cool piece of code:" + [char]10 + "List<Deque<Integer>> pegs = IntStream.range(0, NUM_PEGS).mapToObj(i-> new ArrayDeque<Integer>()).collect(Collectors.toList());" + [char]10
$ws.Range("E13").Value = $note

# Update the current selection to match the saved view state (D13).
$ws.Range("D13").Select()

$wb.Save()
